# Update "Price" (column D) values on the active sheet to reflect the
# latest crypto price snapshot pulled by the scheduled GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "243.44"
    3  = "23.82"
    4  = "5.267"
    5  = "0.05824"
    6  = "6.463"
    7  = "3.339"
    8  = "0.8087"
    9  = "0.9038"
    10 = "0.1379"
    11 = "0.07116"
    12 = "0.03084"
    13 = "0.03031"
    14 = "0.09320"
    17 = "0.04712"
    18 = "0.0006012"
    19 = "0.006187"
    21 = "0.003882"
    22 = "0.00008707"
    23 = "3.562"
    24 = "2.172"
    26 = "0.1304"
    28 = "0.0002340"
    40 = "0.03788"
    41 = "0.006296"
    42 = "0.1050"
    43 = "0.002534"
    44 = "0.006948"
    45 = "0.00005323"
    46 = "0.00000000751"
    47 = "0.5112"
    48 = "0.006944"
    49 = "0.00002102"
    50 = "0.0002002"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    # Keep the cell holding text (matching the source data's inline-string
    # typing) instead of letting Excel auto-convert the numeric-looking
    # literal into a Number.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
